# CDS Input file updates
# - Replace the "ParticipantsTab" query (column B, row 2) with the updated
#   Cypher query (adds diagnosis/file/genomic_info optional matches, re-derives
#   the participant list via a second pass, and sorts the collected sample ids).
# - The new query text is longer (18 lines instead of 12), so the row needs to
#   grow to fit it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newParticipantsQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.platform in ['unknown']`nwith p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN`ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id LIMIT 100"

# Row 2 = ParticipantsTab row; column B holds the "query" text.
$ws.Range("B2").Value = $newParticipantsQuery

# Grow row 2 to fit the longer (18-line) query text.
$ws.Rows.Item(2).RowHeight = 279

Write-Output "Updated ParticipantsTab query and row height."
